{"js": "// Fix \"mathematic helpers\" -> \"mathematics helpers\" (insert a missing \"s\"),\n// and move the Word-managed \"_GoBack\" bookmark (last-edit-location marker)\n// from its old spot at the end of the \"Strict software design patterns\n// application \" bullet to the new edit location, right after \"mathematics\"\n// and before \" helpers\".\n\nconst doc = context.document;\nconst body = doc.body;\n\n// 1. Drop the stale \"_GoBack\" bookmark from its previous location.\nconst oldGoBack = doc.getBookmarkRangeOrNullObject(\"_GoBack\");\noldGoBack.load(\"isNullObject\");\nawait context.sync();\nif (!oldGoBack.isNullObject) {\n  doc.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// 2. Locate the typo and correct it in place.\nconst typoResults = body.search(\"mathematic helpers\", { matchCase: true });\ntypoResults.load(\"text\");\nawait context.sync();\n\nif (typoResults.items.length === 0) {\n  throw new Error('Could not find \"mathematic helpers\" to correct.');\n}\n\nconst typoRange = typoResults.items[0];\ntypoRange.insertText(\"mathematics helpers\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3. Re-insert \"_GoBack\" right after the corrected word \"mathematics\".\nconst fixedResults = body.search(\"mathematics\", { matchCase: true });\nawait context.sync();\n\nif (fixedResults.items.length === 0) {\n  throw new Error('Could not find \"mathematics\" after the correction.');\n}\n\nconst afterFixed = fixedResults.items[0].getRange(\"After\");\nafterFixed.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Fix \"mathematic helpers\" -> \"mathematics helpers\" (insert a missing \"s\"),\n# and move the Word-managed \"_GoBack\" bookmark (last-edit-location marker)\n# from its old spot at the end of the \"Strict software design patterns\n# application \" bullet to the new edit location, right after \"mathematics\"\n# and before \" helpers\".\n\n$d = $word.ActiveDocument\n\n# 1. Drop the stale \"_GoBack\" bookmark from its previous location, if present.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2. Locate the typo and correct it in place.\n$findRange = $d.Content\n$found = $findRange.Find.Execute(\"mathematic helpers\", $false, $false, $false, $false, $false, $true, 1, $false, \"mathematics helpers\", 2)\n\n# 3. Re-insert \"_GoBack\" right after the corrected word \"mathematics\".\n$wordRange = $d.Content\n$wordRange.Find.Execute(\"mathematics\") | Out-Null\n$wordRange.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $wordRange)\n"}
